$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122, shifting the existing rows 122-163 down to 123-164.
$ws.Rows.Item(122).EntireRow.Insert()

# Populate the newly inserted row 122 with the new record's data.
$ws.Range("A122").Value = 7
$ws.Range("B122").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C122").Value = "Ñuble"
$ws.Range("D122").Value = 44468
$ws.Range("E122").Value = 16
$ws.Range("F122").Value = 100112008
$ws.Range("G122").Value = "Coliflor"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 700
$ws.Range("L122").Value = 750
$ws.Range("M122").Value = 725
$ws.Range("N122").Value = "$/unidad"
$ws.Range("O122").Value = "Región del Maule"
$ws.Range("P122").Value = 725
$ws.Range("Q122").Value = 1
$ws.Range("R122").Value = "Hortaliza"
